$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.105.44'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.789.93'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.93'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0690'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0940'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.049.02'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.54'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.791.98'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.102.58'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.00'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.30'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.10'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.49'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.17'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.07%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.67'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.416.04'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.643'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0192'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.66%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.35'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.30%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.04'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.91'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.72%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.38'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +6.86%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0509'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.28%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.07'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.22%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '107.43'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.950.15'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.09%  '
